$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$firstRow = $tbl.Rows.Item(1)
$newRow = $tbl.Rows.Add($firstRow)

$newRow.Cells.Item(1).Range.Text = "DATE"
$newRow.Cells.Item(2).Range.Text = "09-10-22"
